# Trade #58 closed at 2026-02-17 08:47:59 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 58
$summary.Range("B9").Value = 37.93

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 58
$status.Range("G4").Value = 37.93

# --- New trade row data (Trade #58) ---
$rowNum = 59
$tradeNum = 58
$date = "2026-02-17"
$time = "08:47:52"
$strategy = "MarketMaking"
$side = "DOWN"
$entryPrice = 0.33
$exitPrice = 0.33
$status2 = "CLOSED"
$pnlPct = 0
$pnlDollar = 0
$capitalAfter = 99.59999999999999
$entrySlippage = 0
$exitSlippage = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.14

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item($rowNum, 1).Value = $tradeNum
    # Prefix with an apostrophe so the date-shaped string ("2026-02-17") is
    # stored as literal text instead of being auto-converted to a date
    # serial number (matches the rest of the Date column), then reset the
    # cell style so no lingering quote-prefix style index is left behind.
    $ws.Cells.Item($rowNum, 2).Value = "'" + $date
    $ws.Cells.Item($rowNum, 2).Style = "Normal"
    $ws.Cells.Item($rowNum, 3).Value = $time
    $ws.Cells.Item($rowNum, 4).Value = $strategy
    $ws.Cells.Item($rowNum, 5).Value = $side
    $ws.Cells.Item($rowNum, 6).Value = $entryPrice
    $ws.Cells.Item($rowNum, 7).Value = $exitPrice
    $ws.Cells.Item($rowNum, 8).Value = $status2
    $ws.Cells.Item($rowNum, 9).Value = $pnlPct
    $ws.Cells.Item($rowNum, 10).Value = $pnlDollar
    $ws.Cells.Item($rowNum, 11).Value = $capitalAfter
    $ws.Cells.Item($rowNum, 12).Value = $entrySlippage
    $ws.Cells.Item($rowNum, 13).Value = $exitSlippage
    $ws.Cells.Item($rowNum, 14).Value = $confidence
    $ws.Cells.Item($rowNum, 15).Value = $entryReason
    $ws.Cells.Item($rowNum, 16).Value = $exitReason
    $ws.Cells.Item($rowNum, 17).Value = $duration
}
